$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133; existing rows 133.. shift down to 134..
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new price record.
$ws.Range("A133").Value = 10
$ws.Range("B133").Value = "Vega Modelo de Temuco"
$ws.Range("C133").Value = "La Araucanía"
$ws.Range("D133").Value = 44755
$ws.Range("E133").Value = 9
$ws.Range("F133").Value = "Fruta"
$ws.Range("G133").Value = 100104
$ws.Range("H133").Value = "Frutos de pepita"
$ws.Range("I133").Value = 100104003
$ws.Range("J133").Value = "Membrillo"
$ws.Range("K133").Value = "Champion"
$ws.Range("L133").Value = "Primera"
$ws.Range("M133").Value = 65
$ws.Range("N133").Value = 10000
$ws.Range("O133").Value = 10000
$ws.Range("P133").Value = 10000
$ws.Range("Q133").Value = "$/bandeja 18 kilos granel"
$ws.Range("R133").Value = "Región de O'Higgins"
$ws.Range("S133").Value = 556
$ws.Range("T133").Value = 18
